$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (below the header row), shifting all
# existing data rows down by one.
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row with the latest price data entry.
# The leading apostrophe keeps the date as literal text (matching every
# other row in the column) instead of letting Excel coerce it to a date
# serial number.
$ws.Cells.Item(2, 1).Value = "'2026-02-04"
$ws.Cells.Item(2, 2).Value = 783.5
$ws.Cells.Item(2, 3).Value = 1112
$ws.Cells.Item(2, 4).Value = 3610
